$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" -- for the 4 low-priority md files (rows 4-7),
# bump Priority from "low" to "ht" and refresh "Latest Handoff Datetime" (col H)
# on both the zh-cn and de-de status sheets.

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

foreach ($r in 4..7) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-09-06 08:41:32"

    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-09-06 08:41:38"
}
